$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename BR1CNRKPCAL<n> -> BR1CNRKCAL<n> in column A for rows 13-22,
# and drop the thin-border style those cells previously had (cells become
# plain/no-style), matching the committed XML (no more s="1" on A13:C22).
for ($i = 1; $i -le 10; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 1).Value = "BR1CNRKCAL$i"
    $ws.Cells.Item($row, 1).Style = "Normal"
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Style = "Normal"
}

# Update the view: scroll so row 10 is the top-left visible row and select B21
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("B21").Select()
